$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full rebuilt data block (rows 2-9, columns A:T) with the updated TPM-based
# NATMI statistics for the Tnfsf13 -> Tnfrsf13b ligand/receptor pair. The
# original sheet only had 4 data rows (the 4 sending clusters against a
# single "Resolving-Mac" target); the refreshed data expands this to the
# full 4x2 sending x target cluster cross-product (8 rows).
$data = @(
    @("ECs","Tnfsf13","Tnfrsf13b","ECs",3,1,1.096772333333333,3.290317,0.2426185621302128,0.2426185621302128,1,0.3333333333333333,0.4025163333333333,1.207549,0.04753623304791887,0.04753623304791887,0.4414687781147777,3.973219003033,0.01153317251117278,0.01153317251117278),
    @("ECs","Tnfsf13","Tnfrsf13b","Resolving-Mac",3,1,1.096772333333333,3.290317,0.2426185621302128,0.2426185621302128,3,1,8.065052666666666,24.195158,0.9524637669520812,0.9524637669520811,8.845526631676222,79.609739685086,0.2310853896190401,0.23108538961904),
    @("FAPs","Tnfsf13","Tnfrsf13b","ECs",3,1,0.170281,0.510843,0.03766810132102297,0.03766810132102297,1,0.3333333333333333,0.4025163333333333,1.207549,0.04753623304791887,0.04753623304791887,0.06854088375633334,0.6168679538070001,0.001790599642868768,0.001790599642868768),
    @("FAPs","Tnfsf13","Tnfrsf13b","Resolving-Mac",3,1,0.170281,0.510843,0.03766810132102297,0.03766810132102297,3,1,8.065052666666666,24.195158,0.9524637669520812,0.9524637669520811,1.373325233132667,12.359927098194,0.03587750167815421,0.0358775016781542),
    @("MuSCs","Tnfsf13","Tnfrsf13b","ECs",3,1,0.5018676666666667,1.505603,0.1110188577571507,0.1110188577571507,1,0.3333333333333333,0.4025163333333333,1.207549,0.04753623304791887,0.04753623304791887,0.2020099330052222,1.818089397047,0.005277418295057673,0.005277418295057671),
    @("MuSCs","Tnfsf13","Tnfrsf13b","Resolving-Mac",3,1,0.5018676666666667,1.505603,0.1110188577571507,0.1110188577571507,3,1,8.065052666666666,24.195158,0.9524637669520812,0.9524637669520811,4.047589163363778,36.428302470274,0.1057414394620931,0.105741439462093),
    @("Resolving-Mac","Tnfsf13","Tnfrsf13b","ECs",3,1,2.751641333333334,8.254924000000001,0.6086944787916135,0.6086944787916135,1,0.3333333333333333,0.4025163333333333,1.207549,0.04753623304791887,0.04753623304791887,1.107580580141778,9.968225221276001,0.02893504259881965,0.02893504259881965),
    @("Resolving-Mac","Tnfsf13","Tnfrsf13b","Resolving-Mac",3,1,2.751641333333334,8.254924000000001,0.6086944787916135,0.6086944787916135,3,1,8.065052666666666,24.195158,0.9524637669520812,0.9524637669520811,22.19213227311023,199.729190457992,0.5797594361927939,0.5797594361927939)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
